# dca_checks.docx: rename section headings, reword "improved/decreased" totals,
# and sync the per-project change lists with the new totals (DCA refactor).

$d = $word.ActiveDocument

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Text = $newText
}

function Remove-Para($index) {
    $d.Paragraphs.Item($index).Range.Delete()
}

function Insert-ParaBefore($index, $newText) {
    $d.Paragraphs.Item($index).Range.InsertBefore($newText + "`r")
}

# ---------------------------------------------------------------------------
# Section 1: "Departmental DCA Confidence changes this quarter"
# ---------------------------------------------------------------------------
Set-ParaText 2 "Departmental DCA Confidence changes"

Set-ParaText 4 "Improvements"
Set-ParaText 5 "SoT Improved from Amber/Green to Green"
Set-ParaText 6 "A13 Improved from Amber to Amber/Green"
Remove-Para 7                                              # drop "F9 Improved ..." line
Set-ParaText 7 "Columbia Improved from Amber/Green to Green"
Set-ParaText 8 "3 project(s) in total improved"

Set-ParaText 10 "Decreases"
Insert-ParaBefore 11 "A11 Worsened from Amber to Amber/Red"
Set-ParaText 12 "1 project(s) in total have decreased"

Set-ParaText 14 "Missing ratings"

# ---------------------------------------------------------------------------
# Section 2: "SRO Finance confidence Confidence changes this quarter"
# ---------------------------------------------------------------------------
Set-ParaText 17 "SRO Finance confidence Confidence changes"

Set-ParaText 19 "Improvements"
Set-ParaText 20 "Columbia Improved from Amber to Green"
Set-ParaText 21 "1 project(s) in total improved"

Set-ParaText 23 "Decreases"

Set-ParaText 26 "Missing ratings"

# ---------------------------------------------------------------------------
# Section 3: "SRO Benefits RAG Confidence changes this quarter"
# ---------------------------------------------------------------------------
Set-ParaText 29 "SRO Benefits RAG Confidence changes"

Set-ParaText 31 "Improvements"
Set-ParaText 32 "SoT Improved from Amber to Green"
Remove-Para 33                                              # drop "A13 Improved ..." line
Remove-Para 33                                              # drop "F9 Improved ..." line
Remove-Para 33                                              # drop "Columbia Improved ..." line
Set-ParaText 33 "1 project(s) in total improved"

Set-ParaText 35 "Decreases"

Set-ParaText 38 "Missing ratings"
Insert-ParaBefore 39 "A11 Missing"
Set-ParaText 40 "1 project(s) in total are missing a rating"

# ---------------------------------------------------------------------------
# Section 4: "SRO Schedule Confidence Confidence changes this quarter"
# ---------------------------------------------------------------------------
Set-ParaText 42 "SRO Schedule Confidence Confidence changes"

Set-ParaText 44 "Improvements"
Set-ParaText 45 "0 project(s) in total improved"

Set-ParaText 47 "Decreases"

Set-ParaText 50 "Missing ratings"
